$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix duplicated nomenclature values in column J (Informe de RTF) for rows 2-4
$ws.Range("J2").Value = "Plan de Iteración"
$ws.Range("J3").Value = "Informe de Revisión tecnica formal."
$ws.Range("J4").Value = "Informe Final de SQA"

# Update the active selection to J5
$ws.Range("J5").Select()
